# chore: update Sheets via scheduled runner
# Refresh market-price-derived columns (currentAveragePrice[/NQ/HQ],
# LevePriceNQ/HQ, LeveProfitNQ/HQ) for the affected Leve rows across
# all eight job sheets, mirroring the latest scheduled price pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 994
$ws.Range("I86").Value = 988
$ws.Range("K86").Value = 988
$ws.Range("M86").Value = 135
$ws.Range("H89").Value = 994
$ws.Range("I89").Value = 988
$ws.Range("K89").Value = 4940
$ws.Range("M89").Value = 676
$ws.Range("H92").Value = 964.2727
$ws.Range("I92").Value = 964.2727
$ws.Range("K92").Value = 964.2727
$ws.Range("M92").Value = 283.7273
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("M98").ClearContents()
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H137").Value = 3447.4
$ws.Range("I137").Value = 2982.7144
$ws.Range("J137").Value = 4531.6665
$ws.Range("K137").Value = 8948.143199999999
$ws.Range("L137").Value = 13594.9995
$ws.Range("M137").Value = -6398.143199999999
$ws.Range("N137").Value = -18694.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9162.772000000001
$ws.Range("I32").Value = 5093.0586
$ws.Range("K32").Value = 5093.0586
$ws.Range("M32").Value = -4806.0586
$ws.Range("H61").Value = 11674.765
$ws.Range("I61").Value = 8150.4
$ws.Range("K61").Value = 8150.4
$ws.Range("M61").Value = -7938.4
$ws.Range("H74").Value = 5824.25
$ws.Range("I74").Value = 5265.6665
$ws.Range("K74").Value = 5265.6665
$ws.Range("M74").Value = -4391.6665
$ws.Range("H77").Value = 5824.25
$ws.Range("I77").Value = 5265.6665
$ws.Range("K77").Value = 26328.3325
$ws.Range("M77").Value = -21960.3325
$ws.Range("H102").Value = 2113.5
$ws.Range("I102").Value = 2022.5714
$ws.Range("K102").Value = 2022.5714
$ws.Range("M102").Value = -400.5714
$ws.Range("H110").Value = 863.625
$ws.Range("I110").Value = 883
$ws.Range("J110").Value = 805.5
$ws.Range("K110").Value = 883
$ws.Range("L110").Value = 805.5
$ws.Range("M110").Value = 1162
$ws.Range("N110").Value = -4895.5
$ws.Range("H122").Value = 800
$ws.Range("I122").Value = 800
$ws.Range("K122").Value = 2400
$ws.Range("M122").Value = 50
$ws.Range("H136").Value = 11674.765
$ws.Range("I136").Value = 8150.4
$ws.Range("K136").Value = 24451.2
$ws.Range("M136").Value = -21901.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5212.5
$ws.Range("I86").Value = 2750
$ws.Range("J86").Value = 6033.3335
$ws.Range("K86").Value = 2750
$ws.Range("L86").Value = 6033.3335
$ws.Range("M86").Value = -1627
$ws.Range("N86").Value = -8279.333500000001
$ws.Range("H89").Value = 5212.5
$ws.Range("I89").Value = 2750
$ws.Range("J89").Value = 6033.3335
$ws.Range("K89").Value = 13750
$ws.Range("L89").Value = 30166.6675
$ws.Range("M89").Value = -8134
$ws.Range("N89").Value = -41398.6675
$ws.Range("H94").Value = 2062.125
$ws.Range("I94").Value = 1208.6
$ws.Range("J94").Value = 3484.6667
$ws.Range("K94").Value = 1208.6
$ws.Range("L94").Value = 3484.6667
$ws.Range("M94").Value = -757.5999999999999
$ws.Range("N94").Value = -4386.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 25001
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()
$ws.Range("H6").Value = 8060445
$ws.Range("I6").Value = 8060445
$ws.Range("K6").Value = 8060445
$ws.Range("M6").Value = -8060332
$ws.Range("H7").Value = 37.4
$ws.Range("I7").Value = 10.833333
$ws.Range("J7").Value = 77.25
$ws.Range("K7").Value = 10.833333
$ws.Range("L7").Value = 77.25
$ws.Range("M7").Value = 102.166667
$ws.Range("N7").Value = -303.25
$ws.Range("H10").Value = 737.25
$ws.Range("I10").Value = 316.66666
$ws.Range("J10").Value = 1999
$ws.Range("K10").Value = 316.66666
$ws.Range("L10").Value = 1999
$ws.Range("M10").Value = -177.66666
$ws.Range("N10").Value = -2277
$ws.Range("H99").Value = 3524.6667
$ws.Range("J99").Value = 2999.5
$ws.Range("L99").Value = 2999.5
$ws.Range("N99").Value = -5995.5
$ws.Range("H126").Value = 3524.6667
$ws.Range("J126").Value = 2999.5
$ws.Range("L126").Value = 8998.5
$ws.Range("N126").Value = -13938.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H119").Value = 3949.75
$ws.Range("I119").Value = 3949.75
$ws.Range("K119").Value = 11849.25
$ws.Range("M119").Value = -7011.25
$ws.Range("H139").Value = 4777
$ws.Range("I139").Value = 2234.3333
$ws.Range("K139").Value = 6702.999899999999
$ws.Range("M139").Value = -1562.999899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 33950
$ws.Range("I80").Value = 36000
$ws.Range("J80").Value = 32412.5
$ws.Range("K80").Value = 36000
$ws.Range("L80").Value = 32412.5
$ws.Range("M80").Value = -35002
$ws.Range("N80").Value = -34408.5
$ws.Range("H83").Value = 33950
$ws.Range("I83").Value = 36000
$ws.Range("J83").Value = 32412.5
$ws.Range("K83").Value = 180000
$ws.Range("L83").Value = 162062.5
$ws.Range("M83").Value = -175008
$ws.Range("N83").Value = -172046.5
$ws.Range("H97").Value = 1366.1666
$ws.Range("I97").Value = 992.2222
$ws.Range("K97").Value = 992.2222
$ws.Range("M97").Value = -496.2222
$ws.Range("H101").Value = 8399
$ws.Range("J101").Value = 8399
$ws.Range("L101").Value = 8399
$ws.Range("N101").Value = -14889
$ws.Range("H102").Value = 356.3846
$ws.Range("I102").Value = 276.7
$ws.Range("J102").Value = 622
$ws.Range("K102").Value = 276.7
$ws.Range("L102").Value = 622
$ws.Range("M102").Value = 1345.3
$ws.Range("N102").Value = -3866
$ws.Range("H104").Value = 104400
$ws.Range("J104").Value = 104400
$ws.Range("L104").Value = 104400
$ws.Range("N104").Value = -111388
$ws.Range("H132").Value = 1823.1538
$ws.Range("I132").Value = 1336.7273
$ws.Range("K132").Value = 4010.1819
$ws.Range("M132").Value = -1480.1819

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 50
$ws.Range("I2").Value = 50
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 50
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 62
$ws.Range("N2").ClearContents()
$ws.Range("H53").Value = 9000
$ws.Range("I53").Value = 2500
$ws.Range("J53").Value = 15500
$ws.Range("K53").Value = 2500
$ws.Range("L53").Value = 15500
$ws.Range("M53").Value = -1982
$ws.Range("N53").Value = -16536
$ws.Range("H56").Value = 19168.4
$ws.Range("I56").Value = 11746.25
$ws.Range("K56").Value = 11746.25
$ws.Range("M56").Value = -11055.25
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("M82").ClearContents()
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("M85").ClearContents()
$ws.Range("N85").ClearContents()
$ws.Range("H136").Value = 17672.143
$ws.Range("I136").Value = 19984.5
$ws.Range("K136").Value = 59953.5
$ws.Range("M136").Value = -57403.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 3015
$ws.Range("J22").Value = 3015
$ws.Range("L22").Value = 3015
$ws.Range("N22").Value = -3601
$ws.Range("H96").Value = 2323.5557
$ws.Range("I96").Value = 1985.8
$ws.Range("J96").Value = 2745.75
$ws.Range("K96").Value = 1985.8
$ws.Range("L96").Value = 2745.75
$ws.Range("M96").Value = -612.8
$ws.Range("N96").Value = -5491.75
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("H136").Value = 1479.6111
$ws.Range("J136").Value = 1350
$ws.Range("L136").Value = 4050
$ws.Range("N136").Value = -9150
